$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Per-cell value updates (re-generated scores/counts after re-running the
# toy-spam confidence analysis with a minimum support of 5 occurrences).
$changes = @(
    @(3,2,1),
    @(3,3,46),
    @(3,4,46),
    @(3,8,0),
    @(4,2,0.8409090909090909),
    @(4,3,37),
    @(4,4,37),
    @(4,8,7),
    @(4,11,0.8392857142857143),
    @(4,12,47),
    @(4,13,47),
    @(4,17,9),
    @(5,1,"poor"),
    @(5,2,0.7464788732394366),
    @(5,3,53),
    @(5,4,53),
    @(5,8,18),
    @(5,10,"amazing"),
    @(5,11,0.8148148148148148),
    @(5,12,22),
    @(5,13,22),
    @(5,17,5),
    @(6,1,"disappointed"),
    @(6,2,0.7365591397849462),
    @(6,3,137),
    @(6,4,137),
    @(6,8,49),
    @(6,10,"favorite"),
    @(6,11,0.6989247311827957),
    @(6,12,65),
    @(6,13,65),
    @(6,17,28),
    @(7,1,"however"),
    @(7,2,0.71875),
    @(7,3,46),
    @(7,4,46),
    @(7,8,18),
    @(7,10,"classic"),
    @(7,11,0.5094339622641509),
    @(7,17,26),
    @(8,1,"broke"),
    @(8,2,0.7038834951456311),
    @(8,3,145),
    @(8,4,145),
    @(8,8,61),
    @(8,10,"excellent"),
    @(8,11,0.46875),
    @(8,12,30),
    @(8,13,30),
    @(8,17,34),
    @(9,1,"returned"),
    @(9,2,0.6578947368421053),
    @(9,3,25),
    @(9,4,25),
    @(9,8,13),
    @(9,10,"thank"),
    @(9,11,0.4347826086956522),
    @(9,12,30),
    @(9,13,30),
    @(9,17,39),
    @(10,1,"waste"),
    @(10,2,0.6283783783783784),
    @(10,3,93),
    @(10,4,93),
    @(10,8,55),
    @(10,10,"great"),
    @(10,11,0.3459016393442623),
    @(10,12,422),
    @(10,13,422),
    @(10,17,798),
    @(11,1,"junk"),
    @(11,2,0.6),
    @(11,3,33),
    @(11,4,33),
    @(11,8,22),
    @(11,10,"love"),
    @(11,11,0.2998565279770445),
    @(11,12,209),
    @(11,13,209),
    @(11,17,488),
    @(12,10,"loves"),
    @(12,11,0.2489626556016598),
    @(12,12,120),
    @(12,13,120),
    @(12,17,362),
    @(13,1,"instead"),
    @(13,2,0.5833333333333334),
    @(13,3,28),
    @(13,4,28),
    @(13,8,20),
    @(13,11,0.2168674698795181),
    @(13,12,36),
    @(13,13,36),
    @(13,17,130),
    @(14,1,"water"),
    @(14,2,0.5714285714285714),
    @(14,3,24),
    @(14,4,24),
    @(14,8,18),
    @(14,10,"best"),
    @(14,11,0.2166666666666667),
    @(14,12,26),
    @(14,13,26),
    @(14,17,94),
    @(15,1,"small"),
    @(15,2,0.4753623188405797),
    @(15,3,164),
    @(15,4,164),
    @(15,8,181),
    @(15,10,"loved"),
    @(15,11,0.1896024464831804),
    @(15,12,62),
    @(15,13,62),
    @(15,17,265),
    @(16,2,0.462962962962963),
    @(16,3,25),
    @(16,4,25),
    @(16,8,29),
    @(16,10,"friends"),
    @(16,11,0.1375661375661376),
    @(16,12,26),
    @(16,13,26),
    @(16,17,163),
    @(17,1,"less"),
    @(17,2,0.45),
    @(17,3,27),
    @(17,4,27),
    @(17,8,33),
    @(17,10,"christmas"),
    @(17,11,0.09236947791164658),
    @(17,12,23),
    @(17,13,23),
    @(17,17,226),
    @(18,1,"plastic"),
    @(18,2,0.4173228346456693),
    @(18,3,53),
    @(18,4,53),
    @(18,8,74),
    @(18,10,"fun"),
    @(18,11,0.08413672217353199),
    @(18,12,96),
    @(18,13,96),
    @(18,17,1045),
    @(19,1,"broken"),
    @(19,2,0.4096385542168675),
    @(19,3,34),
    @(19,4,34),
    @(19,8,49),
    @(19,10,"game"),
    @(19,11,0.03311688311688311),
    @(19,12,51),
    @(19,13,52),
    @(19,14,0.98),
    @(19,15,0.02000000000000002),
    @(19,16,$true),
    @(19,17,1489),
    @(20,1,"apart"),
    @(20,2,0.3894736842105263),
    @(20,3,37),
    @(20,4,37),
    @(20,8,58),
    @(21,1,"paint"),
    @(21,2,0.3492063492063492),
    @(21,3,22),
    @(21,4,22),
    @(21,8,41),
    @(22,1,"difficult"),
    @(22,2,0.3370786516853932),
    @(22,3,30),
    @(22,4,30),
    @(22,8,59),
    @(23,1,"ok"),
    @(23,2,0.328125),
    @(23,3,42),
    @(23,4,42),
    @(23,8,86),
    @(24,1,"cheap"),
    @(24,2,0.2796208530805687),
    @(24,3,59),
    @(24,4,59),
    @(24,8,152),
    @(25,1,"thought"),
    @(25,2,0.2722772277227723),
    @(25,3,55),
    @(25,4,55),
    @(25,8,147),
    @(26,1,"though"),
    @(26,2,0.2393162393162393),
    @(26,3,28),
    @(26,4,28),
    @(26,8,89),
    @(27,2,0.2244897959183673),
    @(27,3,22),
    @(27,4,22),
    @(27,8,76),
    @(28,1,"item"),
    @(28,2,0.213768115942029),
    @(28,3,59),
    @(28,4,59),
    @(28,8,217),
    @(29,1,"size"),
    @(29,2,0.1958762886597938),
    @(29,3,38),
    @(29,4,38),
    @(29,8,156),
    @(30,1,"could"),
    @(30,2,0.1847133757961783),
    @(30,3,29),
    @(30,4,29),
    @(30,8,128),
    @(31,1,"used"),
    @(31,2,0.1828571428571429),
    @(31,3,32),
    @(31,4,32),
    @(31,8,143),
    @(32,1,"would"),
    @(32,2,0.1810089020771513),
    @(32,3,122),
    @(32,4,122),
    @(32,8,552),
    @(33,1,"money"),
    @(33,2,0.180379746835443),
    @(33,3,57),
    @(33,4,57),
    @(33,8,259),
    @(34,2,0.1582278481012658),
    @(34,3,50),
    @(34,4,50),
    @(34,8,266),
    @(35,1,"better"),
    @(35,2,0.1448598130841121),
    @(35,3,31),
    @(35,4,31),
    @(35,8,183),
    @(36,1,"product"),
    @(36,2,0.1387665198237885),
    @(36,3,63),
    @(36,4,63),
    @(36,8,391),
    @(37,1,"hard"),
    @(37,2,0.13),
    @(37,3,26),
    @(37,4,26),
    @(37,8,174),
    @(38,1,"2"),
    @(38,2,0.1048689138576779),
    @(38,3,28),
    @(38,4,28),
    @(38,8,239),
    @(39,1,"price"),
    @(39,2,0.1005747126436782),
    @(39,3,35),
    @(39,4,35),
    @(39,8,313),
    @(40,1,"use"),
    @(40,2,0.07945205479452055),
    @(40,3,29),
    @(40,4,29),
    @(40,5,0),
    @(40,6,1),
    @(40,7,$false),
    @(40,8,336),
    @(41,1,"like"),
    @(41,2,0.06085526315789474),
    @(41,3,37),
    @(41,4,37),
    @(41,8,571)
)

foreach ($item in $changes) {
    $ws.Cells.Item($item[0], $item[1]).Value = $item[2]
}

# The re-run produced 3 fewer rows in the "negative" word block (41 vs 44),
# so the trailing rows are removed.
$ws.Range("A42:A44").EntireRow.Delete()
